# Applies the "Updated cryptos list" refresh: new Price/Volume(1h) figures
# pulled for this run, plus the Kaspa/ARBITRUM rows (38/39) swapping rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain decimal number (e.g. "2.88") would be
# auto-coerced to a numeric value by Excel on assignment, which would change
# their stored type away from the original text/string cells. Force the
# Text format before writing, then restore the default "Normal" style so no
# visible formatting changes are introduced.
function Set-TextValue($rangeAddress, $value) {
    $cell = $ws.Range($rangeAddress)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '43.204.46'
$ws.Range('E2').Value = '  +0.34%  '
# Row 3
Set-TextValue 'D3' '2.323.86'
# Row 4
$ws.Range('E4').Value = '  +0.00%  '
# Row 5
Set-TextValue 'D5' '302.93'
$ws.Range('E5').Value = '  +0.17%  '
# Row 6
Set-TextValue 'D6' '99.56'
$ws.Range('E6').Value = '  +0.38%  '
# Row 7
$ws.Range('E7').Value = '  +0.17%  '
# Row 8
$ws.Range('E8').Value = '  +0.00%  '
# Row 9
$ws.Range('E9').Value = '  +2.01%  '
# Row 10
Set-TextValue 'D10' '35.97'
$ws.Range('E10').Value = '  +4.66%  '
# Row 11
$ws.Range('E11').Value = '  -0.70%  '
# Row 12
$ws.Range('E12').Value = '  -0.84%  '
# Row 13
Set-TextValue 'D13' '17.66'
$ws.Range('E13').Value = '  -2.48%  '
# Row 14
Set-TextValue 'D14' '6.93'
$ws.Range('E14').Value = '  +1.86%  '
# Row 15
Set-TextValue 'D15' '2.682.97'
$ws.Range('E15').Value = '  +0.86%  '
# Row 16
Set-TextValue 'D16' '2.387.57'
$ws.Range('E16').Value = '  +0.62%  '
# Row 17
Set-TextValue 'D17' '0.798'
$ws.Range('E17').Value = '  -1.41%  '
# Row 18
Set-TextValue 'D18' '43.095.45'
$ws.Range('E18').Value = '  +0.34%  '
# Row 19
Set-TextValue 'D19' '13.22'
$ws.Range('E19').Value = '  +6.61%  '
# Row 20
$ws.Range('E20').Value = '  +2.24%  '
# Row 21
Set-TextValue 'D21' '0.0₃0912'
$ws.Range('E21').Value = '  +0.66%  '
# Row 22
Set-TextValue 'D22' '68.14'
$ws.Range('E22').Value = '  +0.27%  '
# Row 23
Set-TextValue 'D23' '240.40'
$ws.Range('E23').Value = '  +1.58%  '
# Row 24
$ws.Range('E24').Value = '  -2.59%  '
# Row 25
$ws.Range('E25').Value = '  -0.26%  '
# Row 26
$ws.Range('E26').Value = '  -0.14%  '
# Row 27
Set-TextValue 'D27' '25.56'
$ws.Range('E27').Value = '  +3.33%  '
# Row 28
Set-TextValue 'D28' '168.19'
$ws.Range('E28').Value = '  -0.13%  '
# Row 29
Set-TextValue 'D29' '34.31'
$ws.Range('E29').Value = '  +1.48%  '
# Row 30
$ws.Range('E30').Value = '  +0.61%  '
# Row 31
Set-TextValue 'D31' '2.05'
$ws.Range('E31').Value = '  -1.95%  '
# Row 32
Set-TextValue 'D32' '5.18'
$ws.Range('E32').Value = '  +3.03%  '
# Row 33
$ws.Range('E33').Value = '  -0.12%  '
# Row 34
$ws.Range('E34').Value = '  +4.82%  '
# Row 35
Set-TextValue 'D35' '17.72'
$ws.Range('E35').Value = '  +4.44%  '
# Row 36
$ws.Range('E36').Value = '  -0.83%  '
# Row 37
$ws.Range('E37').Value = '  -0.17%  '
# Row 38
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D38' '1.82'
$ws.Range('E38').Value = '  +1.87%  '
# Row 39
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D39' '0.103'
$ws.Range('E39').Value = '  +0.11%  '
# Row 40
$ws.Range('E40').Value = '  -1.44%  '
# Row 41
$ws.Range('E41').Value = '  +0.46%  '
# Row 42
Set-TextValue 'D42' '1.995.96'
$ws.Range('E42').Value = '  -0.19%  '
# Row 43
$ws.Range('E43').Value = '  +1.31%  '
# Row 44
$ws.Range('E44').Value = '  -4.74%  '
# Row 45
Set-TextValue 'D45' '10.11'
$ws.Range('E45').Value = '  +0.61%  '
# Row 46
$ws.Range('E46').Value = '  +0.07%  '
# Row 47
Set-TextValue 'D47' '2.88'
$ws.Range('E47').Value = '  +0.48%  '
# Row 48
Set-TextValue 'D48' '76.68'
$ws.Range('E48').Value = '  +9.14%  '
# Row 49
Set-TextValue 'D49' '55.08'
$ws.Range('E49').Value = '  -0.78%  '
# Row 50
Set-TextValue 'D50' '2.88'
$ws.Range('E50').Value = '  +14.28%  '
# Row 51
Set-TextValue 'D51' '2.548.27'
$ws.Range('E51').Value = '  +0.78%  '
